$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.544.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.253.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.87"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.40%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.628"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0958"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.96%  "

$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.867"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.254.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.445.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("E18").Value = "  +4.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +55.07%  "

$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.54%  "

$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.53%  "

$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("E28").Value = "  +3.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +22.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0827"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.121"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.42%  "

$ws.Range("E36").Value = "  +0.24%  "

$ws.Range("E37").Value = "  +5.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.205"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "108.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "

$ws.Range("E44").Value = "  +3.10%  "

$ws.Range("E45").Value = "  +2.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.40%  "

$ws.Range("E48").Value = "  +1.42%  "

$ws.Range("E49").Value = "  +2.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("E51").Value = "  +1.05%  "
